# ElderRay.Calc.xlsx edit:
#  - rename the "index" column header (A1, and the "testdata" table's
#    first column) to "i"
#  - the index column values become 0-based instead of 1-based
#    (A2:A503 go from 1..502 down to 0..501)
#  - column A is narrowed (bestFit width 6 -> 4 characters)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename header "index" -> "i" (this also updates the ListObject's
# table column name and the shared string table automatically).
$ws.Range("A1").Value = "i"

# Shift the row index values (column A, rows 2-503) down by one so the
# data is numbered starting at 0 instead of 1.
$idxRange = $ws.Range("A2:A503")
$values = $idxRange.Value2
$rowCount = $values.GetLength(0)
for ($r = 1; $r -le $rowCount; $r++) {
    $values[$r, 1] = $values[$r, 1] - 1
}
$idxRange.Value2 = $values

# Narrow column A from its old best-fit width down to 4 characters.
$ws.Columns.Item(1).ColumnWidth = 3.16
